$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 913.17645
$ws.Range("I18").Value = 711.46155
$ws.Range("J18").Value = 1568.75
$ws.Range("K18").Value = 711.46155
$ws.Range("L18").Value = 1568.75
$ws.Range("M18").Value = -427.46155
$ws.Range("N18").Value = -2136.75

# Row 31
$ws.Range("H31").Value = 1037.6666
$ws.Range("I31").Value = 1037.6666
$ws.Range("K31").Value = 3112.9998
$ws.Range("M31").Value = -2882.9998

# Row 40
$ws.Range("H40").Value = 4044.4546
$ws.Range("I40").Value = 2959.8
$ws.Range("J40").Value = 4948.3335
$ws.Range("K40").Value = 2959.8
$ws.Range("L40").Value = 4948.3335
$ws.Range("M40").Value = -2784.8
$ws.Range("N40").Value = -5298.3335

# Row 51
$ws.Range("H51").Value = 150501.28
$ws.Range("I51").Value = 206701.8
$ws.Range("K51").Value = 206701.8
$ws.Range("M51").Value = -206217.8

# Row 128
$ws.Range("H128").Value = 62398
$ws.Range("J128").Value = 62398
$ws.Range("L128").Value = 62398
$ws.Range("N128").Value = -72358

# Row 132
$ws.Range("H132").Value = 1410.5593
$ws.Range("I132").Value = 1563.1765
$ws.Range("K132").Value = 4689.529500000001
$ws.Range("M132").Value = -2159.529500000001

# Row 137
$ws.Range("H137").Value = 3531.75
$ws.Range("I137").Value = 2610.1667
$ws.Range("K137").Value = 7830.500100000001
$ws.Range("M137").Value = -5280.500100000001

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 905.1111
$ws.Range("I2").Value = 939.2941
$ws.Range("J2").Value = 324
$ws.Range("K2").Value = 939.2941
$ws.Range("L2").Value = 324
$ws.Range("M2").Value = -826.2941
$ws.Range("N2").Value = -550

# Row 32
$ws.Range("H32").Value = 6026145
$ws.Range("I32").Value = 6174436
$ws.Range("K32").Value = 6174436
$ws.Range("M32").Value = -6174149

# Row 45
$ws.Range("H45").Value = 1655.9412
$ws.Range("I45").Value = 1506.7273
$ws.Range("K45").Value = 1506.7273
$ws.Range("M45").Value = -1129.7273

# Row 92
$ws.Range("H92").Value = 47329.668
$ws.Range("J92").Value = 47329.668
$ws.Range("L92").Value = 47329.668
$ws.Range("N92").Value = -52321.668

# Row 102
$ws.Range("H102").Value = 8467.299999999999
$ws.Range("I102").Value = 9320.556
$ws.Range("J102").Value = 788
$ws.Range("K102").Value = 9320.556
$ws.Range("L102").Value = 788
$ws.Range("M102").Value = -7698.556
$ws.Range("N102").Value = -4032

# Row 116
$ws.Range("H116").Value = 905.1111
$ws.Range("I116").Value = 939.2941
$ws.Range("J116").Value = 324
$ws.Range("K116").Value = 939.2941
$ws.Range("L116").Value = 324
$ws.Range("M116").Value = 1354.7059
$ws.Range("N116").Value = -4912

# Row 122
$ws.Range("H122").Value = 2253.5312
$ws.Range("I122").Value = 1803.5385
$ws.Range("K122").Value = 5410.6155
$ws.Range("M122").Value = -2960.6155

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 905.1111
$ws.Range("I3").Value = 939.2941
$ws.Range("J3").Value = 324
$ws.Range("K3").Value = 939.2941
$ws.Range("L3").Value = 324
$ws.Range("M3").Value = -825.2941
$ws.Range("N3").Value = -552

# Row 20
$ws.Range("H20").Value = 2332.093
$ws.Range("I20").Value = 2857.875
$ws.Range("J20").Value = 802.5454999999999
$ws.Range("K20").Value = 2857.875
$ws.Range("L20").Value = 802.5454999999999
$ws.Range("M20").Value = -2610.875
$ws.Range("N20").Value = -1296.5455

# Row 86
$ws.Range("H86").Value = 2171.35
$ws.Range("I86").Value = 1848.7858
$ws.Range("K86").Value = 1848.7858
$ws.Range("M86").Value = -725.7858000000001

# Row 89
$ws.Range("H89").Value = 2171.35
$ws.Range("I89").Value = 1848.7858
$ws.Range("K89").Value = 9243.929
$ws.Range("M89").Value = -3627.929

# Row 134
$ws.Range("H134").Value = 669514.4
$ws.Range("I134").Value = 2608.1667
$ws.Range("J134").Value = 3337139.2
$ws.Range("K134").Value = 7824.500100000001
$ws.Range("L134").Value = 10011417.6
$ws.Range("M134").Value = -5289.500100000001
$ws.Range("N134").Value = -10016487.6

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 622136.7
$ws.Range("I31").Value = 10994.5
$ws.Range("J31").Value = 1172164.8
$ws.Range("K31").Value = 10994.5
$ws.Range("L31").Value = 1172164.8
$ws.Range("M31").Value = -10699.5
$ws.Range("N31").Value = -1172754.8

# Row 34
$ws.Range("H34").Value = 622136.7
$ws.Range("I34").Value = 10994.5
$ws.Range("J34").Value = 1172164.8
$ws.Range("K34").Value = 10994.5
$ws.Range("L34").Value = 1172164.8
$ws.Range("M34").Value = -10792.5
$ws.Range("N34").Value = -1172568.8

# Row 62
$ws.Range("H62").Value = 2445.111
$ws.Range("I62").Value = 2388.25
$ws.Range("K62").Value = 2388.25
$ws.Range("M62").Value = -1764.25

# Row 65
$ws.Range("H65").Value = 2445.111
$ws.Range("I65").Value = 2388.25
$ws.Range("K65").Value = 11941.25
$ws.Range("M65").Value = -8821.25

# Row 86
$ws.Range("H86").Value = 5249.1665
$ws.Range("I86").Value = 5249.1665
$ws.Range("K86").Value = 5249.1665
$ws.Range("M86").Value = -4126.1665

# Row 89
$ws.Range("H89").Value = 5249.1665
$ws.Range("I89").Value = 5249.1665
$ws.Range("K89").Value = 26245.8325
$ws.Range("M89").Value = -20629.8325

# Row 105
$ws.Range("H105").Value = 2065.4614
$ws.Range("I105").Value = 2065.4614
$ws.Range("K105").Value = 2065.4614
$ws.Range("M105").Value = -318.4614000000001

# Row 107
$ws.Range("H107").Value = 1751.0667
$ws.Range("I107").Value = 1925.3334
$ws.Range("K107").Value = 1925.3334
$ws.Range("M107").Value = -5.333399999999983

# Row 122
$ws.Range("H122").Value = 978.0909
$ws.Range("I122").Value = 965.9
$ws.Range("K122").Value = 2897.7
$ws.Range("M122").Value = -447.6999999999998

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1599
$ws.Range("J5").Value = 910.7143
$ws.Range("L5").Value = 2732.1429
$ws.Range("N5").Value = -2956.1429

# Row 40
$ws.Range("H40").Value = 172.66667
$ws.Range("I40").Value = 168.07143
$ws.Range("J40").Value = 181.85715
$ws.Range("K40").Value = 672.28572
$ws.Range("L40").Value = 727.4286
$ws.Range("M40").Value = -603.28572
$ws.Range("N40").Value = -865.4286

# Row 68
$ws.Range("H68").Value = 2105.375
$ws.Range("I68").Value = 1490.5
$ws.Range("J68").Value = 3950
$ws.Range("K68").Value = 4471.5
$ws.Range("L68").Value = 11850
$ws.Range("M68").Value = -3660.5
$ws.Range("N68").Value = -13472

# Row 71
$ws.Range("H71").Value = 2105.375
$ws.Range("I71").Value = 1490.5
$ws.Range("J71").Value = 3950
$ws.Range("K71").Value = 13414.5
$ws.Range("L71").Value = 35550
$ws.Range("M71").Value = -9358.5
$ws.Range("N71").Value = -43662

# Row 113
$ws.Range("H113").Value = 2861.5
$ws.Range("I113").Value = 876.5
$ws.Range("J113").Value = 3523.1667
$ws.Range("K113").Value = 2629.5
$ws.Range("L113").Value = 10569.5001
$ws.Range("M113").Value = -459.5
$ws.Range("N113").Value = -14909.5001

# Row 135
$ws.Range("H135").Value = 1599
$ws.Range("J135").Value = 910.7143
$ws.Range("L135").Value = 8196.4287
$ws.Range("N135").Value = -13266.4287

$ws = $wb.Worksheets.Item("GSM")
# Row 59
$ws.Range("H59").Value = 8000
$ws.Range("J59").Value = 8000
$ws.Range("L59").Value = 8000
$ws.Range("N59").Value = -9166

# Row 97
$ws.Range("H97").Value = 1199.0667
$ws.Range("I97").Value = 1429.6666
$ws.Range("J97").Value = 276.66666
$ws.Range("K97").Value = 1429.6666
$ws.Range("L97").Value = 276.66666
$ws.Range("M97").Value = -933.6666
$ws.Range("N97").Value = -1268.66666

# Row 102
$ws.Range("H102").Value = 2771.392
$ws.Range("I102").Value = 2413.721
$ws.Range("K102").Value = 2413.721
$ws.Range("M102").Value = -791.721

# Row 104
$ws.Range("H104").Value = 60917.75
$ws.Range("J104").Value = 60917.75
$ws.Range("L104").Value = 60917.75
$ws.Range("N104").Value = -67905.75

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 6428.5
$ws.Range("I46").Value = 1949
$ws.Range("J46").Value = 12699.8
$ws.Range("K46").Value = 1949
$ws.Range("L46").Value = 12699.8
$ws.Range("M46").Value = -1761
$ws.Range("N46").Value = -13075.8

# Row 60
$ws.Range("H60").Value = 50001
$ws.Range("J60").Value = 50001
$ws.Range("L60").Value = 50001
$ws.Range("N60").Value = -51019

# Row 122
$ws.Range("H122").Value = 4998.3716
$ws.Range("I122").Value = 5153.926
$ws.Range("J122").Value = 4473.375
$ws.Range("K122").Value = 15461.778
$ws.Range("L122").Value = 13420.125
$ws.Range("M122").Value = -13011.778
$ws.Range("N122").Value = -18320.125

# Row 132
$ws.Range("H132").Value = 213412.52
$ws.Range("I132").Value = 6188.5483
$ws.Range("J132").Value = 591291.5
$ws.Range("K132").Value = 18565.6449
$ws.Range("L132").Value = 1773874.5
$ws.Range("M132").Value = -16035.6449
$ws.Range("N132").Value = -1778934.5

$ws = $wb.Worksheets.Item("WVR")
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# Row 51
$ws.Range("H51").Value = 38267
$ws.Range("I51").Value = 35070
$ws.Range("J51").Value = 39332.668
$ws.Range("K51").Value = 35070
$ws.Range("L51").Value = 39332.668
$ws.Range("M51").Value = -34560
$ws.Range("N51").Value = -40352.668
